$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date value from 2023-09-01 (45170) to 2023-09-05 (45174)
# for rows 2-5, matching the target diff.
$ws.Range("C2:C5").Value = 45174
